$wb = $excel.ActiveWorkbook

# --- doSignIn sheet: update fazia's password ---
$wsSignIn = $wb.Worksheets.Item("doSignIn")
$wsSignIn.Range("B2").Value = "Fazia@96us"

# --- createAccount sheet: update last names / emails / passwords ---
$wsCreate = $wb.Worksheets.Item("createAccount")

# row 5 - kate (firstname also changes from "katia " to "kate")
$wsCreate.Range("A5").Value = "kate"
$wsCreate.Range("B5").Value = "katou"
$wsCreate.Range("C5").Value = "katekatrine@gmail.com"
$wsCreate.Range("D5").Value = "Kati1993@!&19"

# row 4 - lysa
$wsCreate.Range("B4").Value = "lulu"
$wsCreate.Range("C4").Value = "lysa.llulu92@gmail.com"
$wsCreate.Range("D4").Value = "Lysa916@%!"

# row 3 - kenza
$wsCreate.Range("B3").Value = "keran"
$wsCreate.Range("C3").Value = "kenza2023@gmail.com"
$wsCreate.Range("D3").Value = "Ken@278!"

# row 2 - linda
$wsCreate.Range("B2").Value = "lady"
$wsCreate.Range("C2").Value = "linda.laydo@gmail.com"
$wsCreate.Range("D2").Value = "Lind@1212#"

# --- make createAccount the active sheet / selected cell, matching workbook's new activeTab ---
$wsCreate.Activate()
[void]$wsCreate.Range("D2").Select()
